$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")
$ws.Cells.Item(1,8).Value = "date"
$ws.Cells.Item(1,9).Value = "legislator_name"
$ws.Cells.Item(1,10).Value = "legislator_id"

$hRange = $ws.Range("H2:H15")
$hRange.NumberFormat = "@"

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r,8).Value = "2012-04-27"
    $ws.Cells.Item($r,9).Value = "江惠贞"
}
